$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "1.00", "0.190", "317.27") keep their exact text representation
# instead of being coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.609.59"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.541.93"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "317.27"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("D6").Value = "94.74"
$ws.Range("E6").Value = "  -3.09%  "
$ws.Range("D7").Value = "0.578"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "36.37"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.114"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "7.55"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").Value = "2.931.64"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "15.66"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").Value = "2.542.95"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "0.866"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "42.654.95"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "13.02"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").Value = "6.61"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "0.0₃0966"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").Value = "70.97"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "254.37"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").Value = "27.34"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  +4.10%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "10.21"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "39.04"
$ws.Range("E30").Value = "  +3.11%  "
$ws.Range("D31").Value = "5.96"
$ws.Range("E31").Value = "  -2.50%  "
$ws.Range("D32").Value = "155.46"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "2.16"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "3.38"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("D35").Value = "19.29"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "0.0787"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("E38").Value = "  -3.44%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "0.119"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "23.88"
$ws.Range("E40").Value = "  -4.07%  "
$ws.Range("D41").Value = "2.33"
$ws.Range("E41").Value = "  +10.75%  "
$ws.Range("D42").Value = "3.83"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0301"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").Value = "2.040.74"
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("D47").Value = "84.53"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "8.92"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "2.788.72"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "74.09"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "0.190"
$ws.Range("E51").Value = "  -0.82%  "

# Restore the default (unstyled) cell style on column D so the
# underlying cell formatting matches the original workbook.
$ws.Range("D2:D51").Style = "Normal"

